$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "-"
$ws.Range("D11").Value = "-"

$ws.Range("B12").Value = "-"
$ws.Range("D12").Value = "-"
$ws.Range("E12").Value = "-"

$ws.Range("D14").Value = "-"
$ws.Range("E14").Value = "-"

$ws.Range("D15").Value = "-"
$ws.Range("E15").Value = "-"

$ws.Range("E16").Value = "-"
